$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record (weekly price observation) as row 277, pushing the
# existing rows 277-339 down to 278-340.
$ws.Rows.Item(277).Insert()

$ws.Range("A277").Value = 4
$ws.Range("B277").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C277").Value = "Los Lagos"
$ws.Range("D277").Value = 44995
$ws.Range("E277").Value = 10
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100108
$ws.Range("H277").Value = "Tropicales y subtropicales"
$ws.Range("I277").Value = 100108002
$ws.Range("J277").Value = "Mango"
$ws.Range("K277").Value = "Sin especificar"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 200
$ws.Range("N277").Value = 8000
$ws.Range("O277").Value = 8500
$ws.Range("P277").Value = 8250
$ws.Range("Q277").Value = "$/bandeja 4 kilos"
$ws.Range("R277").Value = "Perú"
$ws.Range("S277").Value = 2062
$ws.Range("T277").Value = 4
